# PII Example 1 - add Gas columns alongside existing Electric columns,
# and clarify the electric header labels (fix for "no supplier" case).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing Electric headers for clarity -------------------------
$ws.Range("I1").Value = "Electric Choice ID"
$ws.Range("J1").Value = "Electric Rate Code"
$ws.Range("L1").Value = "Electric Usage (kWh)"

# --- Add the three new Gas header cells (M1:O1), matching the existing ---
# --- bold/bordered/centered header style by copying it from K1 -----------
$ws.Range("K1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("N1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("O1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("M1").Value = "Gas Choice ID"
$ws.Range("N1").Value = "Gas Rate Code"
$ws.Range("O1").Value = "Gas Usage (therms)"

# --- Add the (blank) data cells under the new Gas headers for row 2 ------
# A plain empty-string assignment clears/omits the cell entirely, so force
# the cell to materialize as text via a quote-prefixed empty entry, then
# strip the resulting quote-prefix formatting back to the default style so
# it matches the untouched data cells next to it (no bold/border/shading).
$ws.Range("M2").Value = "'"
$ws.Range("M2").Style = "Normal"
$ws.Range("N2").Value = "'"
$ws.Range("N2").Style = "Normal"
$ws.Range("O2").Value = "'"
$ws.Range("O2").Style = "Normal"
